$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-key the two-row sample table (headers in row 1, sample data in row 2).
#    Columns were reshuffled / renamed and a couple of sample values changed.
# ---------------------------------------------------------------------------

# Row 1 - headers
$ws.Range("A1").Value = "Full Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = " Mobile Number"
$ws.Range("D1").Value = "Group"
$ws.Range("E1").Value = "Designation"
$ws.Range("F1").Value = "Gender"
$ws.Range("G1").Value = "Category"
$ws.Range("H1").Value = "Date Of Birth (dd-mm-yyyy)"
$ws.Range("I1").Value = "Mother Tongue"
$ws.Range("J1").Value = "Employee ID"
$ws.Range("K1").Value = "Office Pin Code"
$ws.Range("L1").Value = "External System ID"
$ws.Range("M1").Value = "External System Name"
$ws.Range("N1").Value = "Tags"

# Row 2 - sample data
$ws.Range("A2").Value = "Sahil Chaudhary"
$ws.Range("B2").Value = "sahil11@yopmail.com"
$ws.Range("C2").Value = 7894561230
$ws.Range("D2").Value = "Group A"
$ws.Range("E2").Value = "ACCOUNTANT"
$ws.Range("F2").Value = "Male"
$ws.Range("G2").Value = "OBC"
$ws.Range("H2").Value = "27-07-1998"
$ws.Range("I2").Value = "English"
$ws.Range("J2").Value = "123C"
$ws.Range("K2").Value = 201010
$ws.Range("L2").Value = "USER12345"
$ws.Range("M2").Value = "eHRMSN"
$ws.Range("N2").Value = "Rozgar Mela, Finance"

# ---------------------------------------------------------------------------
# 2. Left-align the " Mobile Number" / "Office Pin Code" columns (header +
#    sample value), matching the new column formatting.
# ---------------------------------------------------------------------------
$ws.Range("C1").HorizontalAlignment = -4131
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("K1").HorizontalAlignment = -4131
$ws.Range("K2").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 3. Move the mailto hyperlink from the name cell (A2) to the email cell (B2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A2").Style = "Normal"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:sahilTest11@yopmail.com") | Out-Null

# ---------------------------------------------------------------------------
# 4. Column widths - reflect the resized / reordered columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.5
$ws.Columns.Item(2).ColumnWidth = 27.666666666666668
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 9.666666666666666
$ws.Columns.Item(5).ColumnWidth = 12.333333333333334
$ws.Columns.Item(6).ColumnWidth = 9
$ws.Columns.Item(7).ColumnWidth = 9.666666666666666
$ws.Columns.Item(8).ColumnWidth = 21.333333333333332
$ws.Columns.Item(9).ColumnWidth = 12.5546875
$ws.Columns.Item(10).ColumnWidth = 10.333333333333334
$ws.Columns.Item(11).ColumnWidth = 12
$ws.Columns.Item(12).ColumnWidth = 16.666666666666668
$ws.Columns.Item(13).ColumnWidth = 20.5546875
$ws.Columns.Item(14).ColumnWidth = 24.333333333333332

# ---------------------------------------------------------------------------
# 5. Selection / scroll position.
# ---------------------------------------------------------------------------
$ws.Range("M2").Select() | Out-Null
